$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.757.65'
$ws.Range('E2').Value = '  +1.25%  '

$ws.Range('D3').Value = '3.924.17'
$ws.Range('E3').Value = '  +1.67%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '468.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +9.79%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.30'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +12.79%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.638'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.38%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('E9').Value = '  +3.27%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.159'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.42%  '

$ws.Range('E11').Value = '  -8.26%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.98'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.32%  '

$ws.Range('D13').Value = '4.556.85'
$ws.Range('E13').Value = '  +1.74%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.45'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.98%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.83'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.01%  '

$ws.Range('D16').Value = '3.912.03'
$ws.Range('E16').Value = '  +0.53%  '

$ws.Range('E17').Value = '  -0.19%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.15'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.13%  '

$ws.Range('E19').Value = '  +8.06%  '

$ws.Range('D20').Value = '68.046.56'
$ws.Range('E20').Value = '  +1.51%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '432.82'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.97%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.85'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.44%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.96%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.82'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.61%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.58'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.50%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '38.24'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.28%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +14.61%  '

$ws.Range('E28').Value = '  +1.48%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.46'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.12%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '735.55'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.48%  '

$ws.Range('E31').Value = '  +11.01%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '13.81'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.38%  '

$ws.Range('E33').Value = '  -0.43%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '43.25'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.90%  '

$ws.Range('E35').Value = '  +7.70%  '

$ws.Range('E36').Value = '  +3.24%  '

$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.51'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.83%  '

$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.26%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0482'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.78%  '

$ws.Range('E40').Value = '  +2.00%  '

$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0693'
$ws.Range('E41').Value = '  -9.13%  '

$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.345'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +10.48%  '

$ws.Range('E43').Value = '  +5.50%  '

$ws.Range('E44').Value = '  +4.55%  '

$ws.Range('E45').Value = '  -0.06%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.55'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +15.34%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.79'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +8.37%  '

$ws.Range('E48').Value = '  +7.19%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.25'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.20%  '

$ws.Range('E50').Value = '  +3.31%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '144.36'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.50%  '
